$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3.0
$ws.Cells.Item(2, 6).Value = 1.0
$ws.Cells.Item(2, 7).Value = 2.11624
$ws.Cells.Item(2, 8).Value = 6.34872
$ws.Cells.Item(2, 9).Value = 0.1897594766532197
$ws.Cells.Item(2, 10).Value = 0.1897594766532197
$ws.Cells.Item(2, 13).Value = 6.322177333333333
$ws.Cells.Item(2, 14).Value = 18.966532
$ws.Cells.Item(2, 15).Value = 0.08271011762055308
$ws.Cells.Item(2, 16).Value = 0.0827101176205531
$ws.Cells.Item(2, 17).Value = 13.37924455989333
$ws.Cells.Item(2, 18).Value = 120.41320103904
$ws.Cells.Item(2, 19).Value = 0.0156950286336024
$ws.Cells.Item(2, 20).Value = 0.0156950286336024

# Row 3
$ws.Cells.Item(3, 5).Value = 3.0
$ws.Cells.Item(3, 6).Value = 1.0
$ws.Cells.Item(3, 7).Value = 2.11624
$ws.Cells.Item(3, 8).Value = 6.34872
$ws.Cells.Item(3, 9).Value = 0.1897594766532197
$ws.Cells.Item(3, 10).Value = 0.1897594766532197
$ws.Cells.Item(3, 15).Value = 0.5401386314560596
$ws.Cells.Item(3, 16).Value = 0.5401386314560597
$ws.Cells.Item(3, 17).Value = 87.3731902988
$ws.Cells.Item(3, 18).Value = 786.3587126892
$ws.Cells.Item(3, 19).Value = 0.1024964240252882
$ws.Cells.Item(3, 20).Value = 0.1024964240252882

# Row 4
$ws.Cells.Item(4, 5).Value = 3.0
$ws.Cells.Item(4, 6).Value = 1.0
$ws.Cells.Item(4, 7).Value = 2.11624
$ws.Cells.Item(4, 8).Value = 6.34872
$ws.Cells.Item(4, 9).Value = 0.1897594766532197
$ws.Cells.Item(4, 10).Value = 0.1897594766532197
$ws.Cells.Item(4, 13).Value = 27.73243066666667
$ws.Cells.Item(4, 14).Value = 83.197292
$ws.Cells.Item(4, 15).Value = 0.3628105447549136
$ws.Cells.Item(4, 16).Value = 0.3628105447549136
$ws.Cells.Item(4, 17).Value = 58.68847907402667
$ws.Cells.Item(4, 18).Value = 528.19631166624
$ws.Cells.Item(4, 19).Value = 0.06884673909696196
$ws.Cells.Item(4, 20).Value = 0.06884673909696196

# Row 5
$ws.Cells.Item(5, 5).Value = 3.0
$ws.Cells.Item(5, 6).Value = 1.0
$ws.Cells.Item(5, 7).Value = 2.11624
$ws.Cells.Item(5, 8).Value = 6.34872
$ws.Cells.Item(5, 9).Value = 0.1897594766532197
$ws.Cells.Item(5, 10).Value = 0.1897594766532197
$ws.Cells.Item(5, 11).Value = 3.0
$ws.Cells.Item(5, 12).Value = 1.0
$ws.Cells.Item(5, 13).Value = 1.096171666666667
$ws.Cells.Item(5, 14).Value = 3.288515
$ws.Cells.Item(5, 15).Value = 0.01434070616847367
$ws.Cells.Item(5, 16).Value = 0.01434070616847367
$ws.Cells.Item(5, 17).Value = 2.319762327866667
$ws.Cells.Item(5, 18).Value = 20.8778609508
$ws.Cells.Item(5, 19).Value = 0.002721284897367162
$ws.Cells.Item(5, 20).Value = 0.002721284897367163

# Row 6
$ws.Cells.Item(6, 9).Value = 0.6160274054778138
$ws.Cells.Item(6, 10).Value = 0.6160274054778138
$ws.Cells.Item(6, 13).Value = 6.322177333333333
$ws.Cells.Item(6, 14).Value = 18.966532
$ws.Cells.Item(6, 15).Value = 0.08271011762055308
$ws.Cells.Item(6, 16).Value = 0.0827101176205531
$ws.Cells.Item(6, 17).Value = 43.4338324433
$ws.Cells.Item(6, 18).Value = 390.9044919897
$ws.Cells.Item(6, 19).Value = 0.05095169916455412
$ws.Cells.Item(6, 20).Value = 0.05095169916455413

# Row 7
$ws.Cells.Item(7, 9).Value = 0.6160274054778138
$ws.Cells.Item(7, 10).Value = 0.6160274054778138
$ws.Cells.Item(7, 15).Value = 0.5401386314560596
$ws.Cells.Item(7, 16).Value = 0.5401386314560597
$ws.Cells.Item(7, 19).Value = 0.3327401997342134
$ws.Cells.Item(7, 20).Value = 0.3327401997342135

# Row 8
$ws.Cells.Item(8, 9).Value = 0.6160274054778138
$ws.Cells.Item(8, 10).Value = 0.6160274054778138
$ws.Cells.Item(8, 13).Value = 27.73243066666667
$ws.Cells.Item(8, 14).Value = 83.197292
$ws.Cells.Item(8, 15).Value = 0.3628105447549136
$ws.Cells.Item(8, 16).Value = 0.3628105447549136
$ws.Cells.Item(8, 17).Value = 190.5238786123
$ws.Cells.Item(8, 18).Value = 1714.7149075107
$ws.Cells.Item(8, 19).Value = 0.2235012385653616
$ws.Cells.Item(8, 20).Value = 0.2235012385653616

# Row 9
$ws.Cells.Item(9, 9).Value = 0.6160274054778138
$ws.Cells.Item(9, 10).Value = 0.6160274054778138
$ws.Cells.Item(9, 11).Value = 3.0
$ws.Cells.Item(9, 12).Value = 1.0
$ws.Cells.Item(9, 13).Value = 1.096171666666667
$ws.Cells.Item(9, 14).Value = 3.288515
$ws.Cells.Item(9, 15).Value = 0.01434070616847367
$ws.Cells.Item(9, 16).Value = 0.01434070616847367
$ws.Cells.Item(9, 17).Value = 7.530781562875001
$ws.Cells.Item(9, 18).Value = 67.777034065875
$ws.Cells.Item(9, 19).Value = 0.008834268013684511
$ws.Cells.Item(9, 20).Value = 0.008834268013684511

# Row 10
$ws.Cells.Item(10, 5).Value = 3.0
$ws.Cells.Item(10, 6).Value = 1.0
$ws.Cells.Item(10, 7).Value = 1.793503666666667
$ws.Cells.Item(10, 8).Value = 5.380511
$ws.Cells.Item(10, 9).Value = 0.1608202836929164
$ws.Cells.Item(10, 10).Value = 0.1608202836929164
$ws.Cells.Item(10, 13).Value = 6.322177333333333
$ws.Cells.Item(10, 14).Value = 18.966532
$ws.Cells.Item(10, 15).Value = 0.08271011762055308
$ws.Cells.Item(10, 16).Value = 0.0827101176205531
$ws.Cells.Item(10, 17).Value = 11.33884822865022
$ws.Cells.Item(10, 18).Value = 102.049634057852
$ws.Cells.Item(10, 19).Value = 0.01330146458001183
$ws.Cells.Item(10, 20).Value = 0.01330146458001183

# Row 11
$ws.Cells.Item(11, 5).Value = 3.0
$ws.Cells.Item(11, 6).Value = 1.0
$ws.Cells.Item(11, 7).Value = 1.793503666666667
$ws.Cells.Item(11, 8).Value = 5.380511
$ws.Cells.Item(11, 9).Value = 0.1608202836929164
$ws.Cells.Item(11, 10).Value = 0.1608202836929164
$ws.Cells.Item(11, 15).Value = 0.5401386314560596
$ws.Cells.Item(11, 16).Value = 0.5401386314560597
$ws.Cells.Item(11, 17).Value = 74.04837691814834
$ws.Cells.Item(11, 18).Value = 666.4353922633351
$ws.Cells.Item(11, 19).Value = 0.0868652479442671
$ws.Cells.Item(11, 20).Value = 0.08686524794426712

# Row 12
$ws.Cells.Item(12, 5).Value = 3.0
$ws.Cells.Item(12, 6).Value = 1.0
$ws.Cells.Item(12, 7).Value = 1.793503666666667
$ws.Cells.Item(12, 8).Value = 5.380511
$ws.Cells.Item(12, 9).Value = 0.1608202836929164
$ws.Cells.Item(12, 10).Value = 0.1608202836929164
$ws.Cells.Item(12, 13).Value = 27.73243066666667
$ws.Cells.Item(12, 14).Value = 83.197292
$ws.Cells.Item(12, 15).Value = 0.3628105447549136
$ws.Cells.Item(12, 16).Value = 0.3628105447549136
$ws.Cells.Item(12, 17).Value = 49.73821608624579
$ws.Cells.Item(12, 18).Value = 447.643944776212
$ws.Cells.Item(12, 19).Value = 0.05834729473426674
$ws.Cells.Item(12, 20).Value = 0.05834729473426674

# Row 13
$ws.Cells.Item(13, 5).Value = 3.0
$ws.Cells.Item(13, 6).Value = 1.0
$ws.Cells.Item(13, 7).Value = 1.793503666666667
$ws.Cells.Item(13, 8).Value = 5.380511
$ws.Cells.Item(13, 9).Value = 0.1608202836929164
$ws.Cells.Item(13, 10).Value = 0.1608202836929164
$ws.Cells.Item(13, 11).Value = 3.0
$ws.Cells.Item(13, 12).Value = 1.0
$ws.Cells.Item(13, 13).Value = 1.096171666666667
$ws.Cells.Item(13, 14).Value = 3.288515
$ws.Cells.Item(13, 15).Value = 0.01434070616847367
$ws.Cells.Item(13, 16).Value = 0.01434070616847367
$ws.Cells.Item(13, 17).Value = 1.965987903462778
$ws.Cells.Item(13, 18).Value = 17.693891131165
$ws.Cells.Item(13, 19).Value = 0.002306276434370691
$ws.Cells.Item(13, 20).Value = 0.002306276434370691

# Row 14
$ws.Cells.Item(14, 5).Value = 2.0
$ws.Cells.Item(14, 6).Value = 0.6666666666666666
$ws.Cells.Item(14, 7).Value = 0.3724043333333333
$ws.Cells.Item(14, 8).Value = 1.117213
$ws.Cells.Item(14, 9).Value = 0.03339283417605023
$ws.Cells.Item(14, 10).Value = 0.03339283417605023
$ws.Cells.Item(14, 13).Value = 6.322177333333333
$ws.Cells.Item(14, 14).Value = 18.966532
$ws.Cells.Item(14, 15).Value = 0.08271011762055308
$ws.Cells.Item(14, 16).Value = 0.0827101176205531
$ws.Cells.Item(14, 17).Value = 2.354406235035111
$ws.Cells.Item(14, 18).Value = 21.189656115316
$ws.Cells.Item(14, 19).Value = 0.002761925242384739
$ws.Cells.Item(14, 20).Value = 0.00276192524238474

# Row 15
$ws.Cells.Item(15, 5).Value = 2.0
$ws.Cells.Item(15, 6).Value = 0.6666666666666666
$ws.Cells.Item(15, 7).Value = 0.3724043333333333
$ws.Cells.Item(15, 8).Value = 1.117213
$ws.Cells.Item(15, 9).Value = 0.03339283417605023
$ws.Cells.Item(15, 10).Value = 0.03339283417605023
$ws.Cells.Item(15, 15).Value = 0.5401386314560596
$ws.Cells.Item(15, 16).Value = 0.5401386314560597
$ws.Cells.Item(15, 17).Value = 15.37545584831167
$ws.Cells.Item(15, 18).Value = 138.379102634805
$ws.Cells.Item(15, 19).Value = 0.0180367597522909
$ws.Cells.Item(15, 20).Value = 0.01803675975229091

# Row 16
$ws.Cells.Item(16, 5).Value = 2.0
$ws.Cells.Item(16, 6).Value = 0.6666666666666666
$ws.Cells.Item(16, 7).Value = 0.3724043333333333
$ws.Cells.Item(16, 8).Value = 1.117213
$ws.Cells.Item(16, 9).Value = 0.03339283417605023
$ws.Cells.Item(16, 10).Value = 0.03339283417605023
$ws.Cells.Item(16, 13).Value = 27.73243066666667
$ws.Cells.Item(16, 14).Value = 83.197292
$ws.Cells.Item(16, 15).Value = 0.3628105447549136
$ws.Cells.Item(16, 16).Value = 0.3628105447549136
$ws.Cells.Item(16, 17).Value = 10.32767735413289
$ws.Cells.Item(16, 18).Value = 92.949096187196
$ws.Cells.Item(16, 19).Value = 0.01211527235832328
$ws.Cells.Item(16, 20).Value = 0.01211527235832328

# Row 17
$ws.Cells.Item(17, 5).Value = 2.0
$ws.Cells.Item(17, 6).Value = 0.6666666666666666
$ws.Cells.Item(17, 7).Value = 0.3724043333333333
$ws.Cells.Item(17, 8).Value = 1.117213
$ws.Cells.Item(17, 9).Value = 0.03339283417605023
$ws.Cells.Item(17, 10).Value = 0.03339283417605023
$ws.Cells.Item(17, 11).Value = 3.0
$ws.Cells.Item(17, 12).Value = 1.0
$ws.Cells.Item(17, 13).Value = 1.096171666666667
$ws.Cells.Item(17, 14).Value = 3.288515
$ws.Cells.Item(17, 15).Value = 0.01434070616847367
$ws.Cells.Item(17, 16).Value = 0.01434070616847367
$ws.Cells.Item(17, 17).Value = 0.4082190787438889
$ws.Cells.Item(17, 18).Value = 3.673971708695
$ws.Cells.Item(17, 19).Value = 0.0004788768230513017
$ws.Cells.Item(17, 20).Value = 0.0004788768230513017
